$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 54155
$ws.Cells.Item(2,2).Value = "Manuela Moreira"
$ws.Cells.Item(2,3).Value = "Recursos Humanos"
$ws.Cells.Item(2,4).Value = "Problemas pessoais"
$ws.Cells.Item(2,5).Value = 8
$ws.Cells.Item(2,6).Value = 45091
$ws.Cells.Item(2,7).Value = 4828.61

# Row 3
$ws.Cells.Item(3,1).Value = 75368
$ws.Cells.Item(3,2).Value = "Heloísa Azevedo"
$ws.Cells.Item(3,3).Value = "TI"
$ws.Cells.Item(3,4).Value = "Viagem de negócios"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 45092
$ws.Cells.Item(3,7).Value = 4839.21

# Row 4
$ws.Cells.Item(4,1).Value = 60405
$ws.Cells.Item(4,2).Value = "Maria Ramos"
$ws.Cells.Item(4,3).Value = "Jurídico"
$ws.Cells.Item(4,4).Value = "Doença"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 45092
$ws.Cells.Item(4,7).Value = 11081.9

# Row 5
$ws.Cells.Item(5,1).Value = 36153
$ws.Cells.Item(5,2).Value = "Stephany Mendes"
$ws.Cells.Item(5,3).Value = "Jurídico"
$ws.Cells.Item(5,4).Value = "Viagem de negócios"
$ws.Cells.Item(5,5).Value = 6
$ws.Cells.Item(5,6).Value = 45089
$ws.Cells.Item(5,7).Value = 6755.53

# Row 6
$ws.Cells.Item(6,1).Value = 68439
$ws.Cells.Item(6,2).Value = "Pedro Monteiro"
$ws.Cells.Item(6,3).Value = "Recursos Humanos"
$ws.Cells.Item(6,4).Value = "Consulta médica"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 45082
$ws.Cells.Item(6,7).Value = 4265.42

# Row 7
$ws.Cells.Item(7,1).Value = 80062
$ws.Cells.Item(7,2).Value = "Davi Luiz Nunes"
$ws.Cells.Item(7,3).Value = "TI"
$ws.Cells.Item(7,4).Value = "Outros"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 45097
$ws.Cells.Item(7,7).Value = 7469.52

# Row 8
$ws.Cells.Item(8,1).Value = 40504
$ws.Cells.Item(8,2).Value = "Sr. Emanuel da Mata"
$ws.Cells.Item(8,3).Value = "P&D"
$ws.Cells.Item(8,4).Value = "Consulta médica"
$ws.Cells.Item(8,5).Value = 7
$ws.Cells.Item(8,6).Value = 45097
$ws.Cells.Item(8,7).Value = 5520.91

# Row 9
$ws.Cells.Item(9,1).Value = 30021
$ws.Cells.Item(9,2).Value = "Raul Fernandes"
$ws.Cells.Item(9,3).Value = "Marketing"
$ws.Cells.Item(9,4).Value = "Problemas pessoais"
$ws.Cells.Item(9,5).Value = 8
$ws.Cells.Item(9,6).Value = 45080
$ws.Cells.Item(9,7).Value = 6682.39

# Row 10
$ws.Cells.Item(10,1).Value = 67036
$ws.Cells.Item(10,2).Value = "Lorena Moraes"
$ws.Cells.Item(10,3).Value = "TI"
$ws.Cells.Item(10,4).Value = "Outros"
$ws.Cells.Item(10,5).Value = 4
$ws.Cells.Item(10,6).Value = 45105
$ws.Cells.Item(10,7).Value = 4030.95

# Row 11
$ws.Cells.Item(11,1).Value = 95015
$ws.Cells.Item(11,2).Value = "Sophie Gomes"
$ws.Cells.Item(11,3).Value = "Operações"
$ws.Cells.Item(11,4).Value = "Consulta médica"
$ws.Cells.Item(11,5).Value = 4
$ws.Cells.Item(11,6).Value = 45096
$ws.Cells.Item(11,7).Value = 8998.370000000001
